$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("A2").Value = 2071139
$ws.Range("B2").Value = 26517
$ws.Range("C2").Value = -5

# Update the active selection to D6
$ws.Range("D6").Select()
